# Apply the cryptos.xlsx price/volume refresh described by the commit
# "Updated cryptos list on Mon Sep 18 23:56:18 UTC 2023 with GitHub Actions".
#
# Rows 2-49: price (column D) and/or 1h volume % change (column E) updates.
# Rows 50-51: EnergySwap and Algorand swapped positions (name/link/price/volume).
#
# Values in column D that look like plain numbers (e.g. "4.12") are written
# with a leading apostrophe so Excel stores them as text, matching the
# workbook's existing convention of keeping the Price column as text
# (note values such as "26.819.94" use "." as a thousands separator, so
# they cannot be real numeric cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.819.94"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "1.641.23"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").Value = "'216.47"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("E6").Value = "  +1.60%  "
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("E8").Value = "  +1.36%  "
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("E10").Value = "  +3.94%  "
$ws.Range("D11").Value = "'0.0844"
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").Value = "1.870.93"
$ws.Range("D13").Value = "1.646.46"
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").Value = "'4.12"
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").Value = "'66.02"
$ws.Range("E16").Value = "  +3.11%  "
$ws.Range("D17").Value = "26.846.24"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").Value = "'218.62"
$ws.Range("E19").Value = "  +3.38%  "
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("E21").Value = "  +1.06%  "
$ws.Range("E22").Value = "  +6.33%  "
$ws.Range("E23").Value = "  +2.39%  "
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("D25").Value = "'146.14"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("E26").Value = "  -0.55%  "
$ws.Range("D27").Value = "'7.42"
$ws.Range("E27").Value = "  +6.04%  "
$ws.Range("E28").Value = "  +1.37%  "
$ws.Range("E29").Value = "  +1.48%  "
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("E31").Value = "  -0.29%  "
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("E33").Value = "  +1.39%  "
$ws.Range("E34").Value = "  +1.68%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").Value = "1.240.83"
$ws.Range("E36").Value = "  -1.72%  "
$ws.Range("E37").Value = "  +0.99%  "
$ws.Range("D38").Value = "'0.535"
$ws.Range("E38").Value = "  +1.87%  "
$ws.Range("D39").Value = "'0.826"
$ws.Range("E39").Value = "  +2.83%  "
$ws.Range("E40").Value = "  -0.43%  "
$ws.Range("D41").Value = "'0.805"
$ws.Range("E41").Value = "  +0.61%  "
$ws.Range("D42").Value = "'5.35"
$ws.Range("E42").Value = "  +1.67%  "
$ws.Range("D43").Value = "1.782.34"
$ws.Range("E43").Value = "  +0.63%  "
$ws.Range("E44").Value = "  -2.85%  "
$ws.Range("D45").Value = "'60.72"
$ws.Range("E45").Value = "  +1.36%  "
$ws.Range("D46").Value = "'91.24"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("E48").Value = "  +7.81%  "
$ws.Range("E49").Value = "  -0.52%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.61"
$ws.Range("E50").Value = "  +2.46%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0969"
$ws.Range("E51").Value = "  +1.12%  "

